$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Reorder / retitle the "Violation" helper rows in columns J:K and turn
#    the old placeholder "Risk Score" row into a real "Total Risk Score"
#    row driven by an (intentionally unresolved-name) array formula.
#    Column J holds the column name, column K holds a description / the
#    actual formula text (as literal text, except for row 10 which becomes
#    a genuine formula).
# ---------------------------------------------------------------------------

# Row 3: Password Age (days)  -- content unchanged, only formatting/order shifts later
$ws.Range("J3").Value = "Password Age (days)"
$ws.Range("K3").Value = "'=TODAY() - [Password Last Changed]"

# Row 4: Login Age (days) -- content unchanged
$ws.Range("J4").Value = "Login Age (days)"
$ws.Range("K4").Value = "'=TODAY() - [Last Login]"

# Row 5: now "Violation: MFA" (moved up from old row 6)
$ws.Range("J5").Value = "Violation: MFA"
$ws.Range("K5").Value = "'=IF([MFA Enabled]=" + [char]34 + "N" + [char]34 + ", " + [char]34 + "YES" + [char]34 + ", " + [char]34 + "NO" + [char]34 + ")"

# Row 6: now "Violation: Password Age" (moved down from old row 5)
$ws.Range("J6").Value = "Violation: Password Age"
$ws.Range("K6").Value = "'=IF([Password Age]>90, " + [char]34 + "YES" + [char]34 + ", " + [char]34 + "NO" + [char]34 + ")"

# Row 7: Violation: Inactive -- content unchanged
$ws.Range("J7").Value = "Violation: Inactive"
$ws.Range("K7").Value = "'=IF([Login Age]>60, " + [char]34 + "YES" + [char]34 + ", " + [char]34 + "NO" + [char]34 + ")"

# Row 8: Violation: Terminated -- content unchanged
$ws.Range("J8").Value = "Violation: Terminated"
$ws.Range("K8").Value = "'=IF([Termination Date]<>" + [char]34 + [char]34 + ", " + [char]34 + "YES" + [char]34 + ", " + [char]34 + "NO" + [char]34 + ")"

# Row 9: Violation: Admin Review -- content unchanged
$ws.Range("J9").Value = "Violation: Admin Review"
$ws.Range("K9").Value = "'=IF(AND([Admin Access]=" + [char]34 + "Y" + [char]34 + ", [Login Age]>30), " + [char]34 + "YES" + [char]34 + ", " + [char]34 + "NO" + [char]34 + ")"

# Row 10: "Risk Score" -> "Total Risk Score", now backed by a real (array)
# formula referencing undefined pseudo-names, which resolves to #NAME?
# -- this mirrors the author's intent of wiring up an actual weighted
# scoring formula while keeping it an unresolved placeholder.
$ws.Range("J10").Value = "Total Risk Score"
$ws.Range("K10").FormulaArray = "=(MFA * 3) + (Password Age * 2) + (Inactive * 1) + (Terminated * 3) + (Admin Review * 2)"

# ---------------------------------------------------------------------------
# 2) New column names (J) become bold; K keeps its existing (Arial Unicode
#    MS) font -- matches the author swapping which style slot is bold.
# ---------------------------------------------------------------------------
$ws.Range("J3:J10").Font.Bold = $true
$ws.Range("K3:K10").Font.Bold = $false

# ---------------------------------------------------------------------------
# 3) Row heights: several helper rows shrink from large wrapped heights to
#    smaller, more precise auto-fit heights; a couple reset to the sheet
#    default entirely.
# ---------------------------------------------------------------------------
$ws.Rows(2).RowHeight = 14.5
$ws.Rows(3).RowHeight = 25
$ws.Rows(4).RowHeight = 14.5
$ws.Rows(5).RowHeight = 14.5
$ws.Rows(6).RowHeight = 29
$ws.Rows(7).RowHeight = 25
$ws.Rows(8).RowHeight = 29
$ws.Rows(9).RowHeight = 37.5
$ws.Rows(10).RowHeight = 14.5

# ---------------------------------------------------------------------------
# 4) New helper column G gets an explicit (best-fit-like) width, mirroring
#    the new data that was fit into that column.
# ---------------------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 9.17

# ---------------------------------------------------------------------------
# 5) View state: zoom level and active selection moved.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 68
$ws.Range("F7").Select()
